$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header column C ("metier_level_6" -> "metier_level_6_new")
$ws.Range("C1").Value = "metier_level_6_new"

# Remove the "LLS_DEF_0_0_0" row entirely (was row 6); remaining rows shift up by one
$ws.Rows.Item(6).Delete()

# --- After the deletion, row numbers have shifted up by 1 for everything below row 6 ---
# Update n_count / KG_sum / EUR_sum for GNS_DEF_110-156_0_0 (row 2)
$ws.Cells.Item(2, 4).Value = 4
$ws.Cells.Item(2, 5).Value = 1600
$ws.Cells.Item(2, 6).Value = 0

# Update n_count / KG_sum / EUR_sum for GNS_FWS_>0_0_0 (row 4)
$ws.Cells.Item(4, 4).Value = 7
$ws.Cells.Item(4, 5).Value = 1080
$ws.Cells.Item(4, 6).Value = 1960

# Update n_count / KG_sum for PTB_DEF_105-115_1_120 (row 9, was row 10 before delete)
$ws.Cells.Item(9, 4).Value = 3
$ws.Cells.Item(9, 5).Value = 900

# Last row (row 11, was row 12 before delete): metier changes from
# "OTB_DWS_100-119_0_0" to "OTB_DEF_105-115_1_120" (reusing the existing string)
$ws.Range("C11").Value = "OTB_DEF_105-115_1_120"
